$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.10"
$ws.Range("D3").Value = "'22.78"
$ws.Range("D4").Value = "'6.193"
$ws.Range("D5").Value = "'0.06112"
$ws.Range("D6").Value = "'6.740"
$ws.Range("D7").Value = "'3.503"
$ws.Range("D8").Value = "'1.357"
$ws.Range("D9").Value = "'0.7990"
$ws.Range("D10").Value = "'0.1580"
$ws.Range("D11").Value = "'0.08062"
$ws.Range("D12").Value = "'0.03307"
$ws.Range("D13").Value = "'0.03036"
$ws.Range("D14").Value = "'0.09299"
$ws.Range("D15").Value = "'3.914"
$ws.Range("D16").Value = "'0.001698"
$ws.Range("D17").Value = "'0.04832"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006144"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006208"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001100"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.003401"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.689"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.260"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D25").Value = "'0.3358"
$ws.Range("D26").Value = "'0.1227"
$ws.Range("D27").Value = "'0.0006169"
$ws.Range("D40").Value = "'0.04592"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1123"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").Value = "'0.003131"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003414"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "'0.01020"
$ws.Range("D46").Value = "'0.00006021"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D48").Value = "'0.7505"
$ws.Range("D49").Value = "'0.1169"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D51").Value = "'0.01011"
